# Weekly update: insert a new price record for "Piña" at Feria Lagunitas de
# Puerto Montt, pushing the existing history (rows 237:250) down by one row
# (to 238:251) and populating the new row 237 with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("237").Insert()

$ws.Range("A237").Value = 4
$ws.Range("B237").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C237").Value = "Los Lagos"
$ws.Range("D237").Value = 44714
$ws.Range("E237").Value = 10
$ws.Range("F237").Value = "Fruta"
$ws.Range("G237").Value = 100108
$ws.Range("H237").Value = "Tropicales y subtropicales"
$ws.Range("I237").Value = 100108005
$ws.Range("J237").Value = "Piña"
$ws.Range("K237").Value = "Caramelo"
$ws.Range("L237").Value = "Tercera"
$ws.Range("M237").Value = 100
$ws.Range("N237").Value = 20000
$ws.Range("O237").Value = 21000
$ws.Range("P237").Value = 20500
$ws.Range("Q237").Value = "$/caja 16 unidades"
$ws.Range("R237").Value = "Ecuador"
$ws.Range("S237").Value = 1281
$ws.Range("T237").Value = 16
